$wb = $excel.ActiveWorkbook

# --- "grilla de pruebas" sheet: flip position to BUY and raise leverage to 30 so
#     the "ataque" (attack/exposure) calculation can extend further down the grid ---
$ws = $wb.Worksheets.Item("grilla de pruebas")

$ws.Range("B1").Value = "BUY"
$ws.Range("B3").Value = 1932.25
$ws.Range("F3").Value = 5.992
$ws.Range("B5").Value = 2
$ws.Range("B9").Value = 5

# Rows 11 and 12 were empty; continue the same E/F progression used by rows 4-10
$ws.Range("E11").Formula = "=E10*(1+`$B`$6/100)"
$ws.Range("F11").Formula = "=IF(`$B`$1=""BUY"",F10*(1-`$B`$5/100),F10*(1+`$B`$5/100))"
$ws.Range("E12").Formula = "=E11*(1+`$B`$6/100)"
$ws.Range("F12").Formula = "=IF(`$B`$1=""BUY"",F11*(1-`$B`$5/100),F11*(1+`$B`$5/100))"

# "ataque" now sums through row 16 instead of stopping at row 9
$ws.Range("E17").Formula = "=SUM(E3:E16)*3"

$ws.Range("E9").Select()

# --- "GRILLA" sheet keeps a text note mirroring that formula; update it to match ---
$wsGrilla = $wb.Worksheets.Item("GRILLA")
$wsGrilla.Range("C17").Value = "ataque seria =SUM(E3:E16)*3"
